$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "gas" for the gas-type generator rows
$ws.Range("B8").Value = "gas"
$ws.Range("B9").Value = "gas"
$ws.Range("B10").Value = "gas"
$ws.Range("B23").Value = "gas"
$ws.Range("B24").Value = "gas"

# Fill in "steam" for the steam-type generator rows
$ws.Range("B7").Value = "steam"
$ws.Range("B12").Value = "steam"
$ws.Range("B17").Value = "steam"
$ws.Range("B19").Value = "steam"

# Fill in "wind" for the wind-type generator rows
$ws.Range("B13").Value = "wind"
$ws.Range("B14").Value = "wind"
$ws.Range("B15").Value = "wind"
$ws.Range("B16").Value = "wind"

# Header cell for the new column
$ws.Range("B4").Value = "LTDtype?"

# "unknown" row gets a dash placeholder
$ws.Range("B6").Value = "-"

# Hydro row reuses the existing "hydro" shared string
$ws.Range("B22").Value = "hydro"

# Restore the selection to match the editing session that produced this workbook
$ws.Range("D8:D10").Select()
